# Auto-generated edit script: update crypto price/volume figures
# per the scraped GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.448.99"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "2.280.72"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D5").Value = "'303.60"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").Value = "'95.07"
$ws.Range("E6").Value = "  -3.16%  "
$ws.Range("D7").Value = "'0.500"
$ws.Range("E7").Value = "  -3.36%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.491"
$ws.Range("E9").Value = "  -3.89%  "
$ws.Range("D10").Value = "'34.80"
$ws.Range("E10").Value = "  -4.40%  "
$ws.Range("D11").Value = "'0.0778"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").Value = "'0.119"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "'18.03"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "'6.65"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "2.632.22"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "2.277.49"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "'0.768"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "42.344.54"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "'12.65"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("D21").Value = "'5.96"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("D22").Value = "'66.99"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").Value = "'235.38"
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("D26").Value = "'2.38"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("D27").Value = "'24.61"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.38"
$ws.Range("E28").Value = "  +17.15%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'166.67"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'8.92"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("D31").Value = "'32.15"
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "'4.92"
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("D34").Value = "'17.53"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("D35").Value = "'4.56"
$ws.Range("E35").Value = "  -4.04%  "
$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("D37").Value = "'0.0678"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").Value = "'0.108"
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("D41").Value = "'2.64"
$ws.Range("E41").Value = "  -4.85%  "
$ws.Range("D42").Value = "1.986.70"
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").Value = "'0.0274"
$ws.Range("E43").Value = "  -4.03%  "
$ws.Range("D44").Value = "'9.98"
$ws.Range("E44").Value = "  -2.70%  "
$ws.Range("D45").Value = "'17.57"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Value = "'2.03"
$ws.Range("E46").Value = "  -7.09%  "
$ws.Range("D47").Value = "'2.73"
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").Value = "'2.91"
$ws.Range("E48").Value = "  +9.69%  "
$ws.Range("D49").Value = "'53.28"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "2.499.13"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "'70.40"
$ws.Range("E51").Value = "  -3.02%  "
